$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.9999999957107899,
    0.99824863123661856,
    0.99243400044790953,
    0.99489028011542624,
    0.98417524324512362,
    0.95985470066373568,
    0.95450763659462134,
    0.94716225630673567,
    0.93852241359592192,
    0.93079562193781129,
    0.92968017570723616,
    0.92781399039010393,
    0.92098197335054444,
    0.91882682565080276,
    0.91623536709455289,
    0.91372889674270907,
    0.91002097444916008,
    0.90891205147383358,
    0.99719539540526381,
    0.99007830706453781,
    0.98867980472420691,
    0.98741529712073262,
    0.98539971215768385,
    0.9723793597748821,
    0.96592244504731883,
    0.95748976550954656,
    0.95444159799517214,
    0.94343533457226703,
    0.93592154613602396,
    0.9293523813402107,
    0.92169942893772239,
    0.92002011450889754,
    0.91950009895148976
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
